$wb = $excel.ActiveWorkbook

# Rename Sheet1 to "FTP Parameterization"
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "FTP Parameterization"

# Populate the data column
$ws.Range("A1").Value = "File Name"
$ws.Range("A2").Value = "a11"
$ws.Range("A3").Value = "b12"
$ws.Range("A4").Value = "c13"
$ws.Range("A5").Value = "d14"

# Set column A width to match bestFit width (stored width = 10)
$ws.Columns.Item(1).ColumnWidth = 9.166666666666666
